# Weekly data refresh: insert a new daily price record for Coliflor at
# "Vega Monumental Concepción" ahead of the existing row 218, pushing all
# subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 218; everything from the old row 218
# through the old row 259 shifts down to rows 219-260.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new record.
$ws.Range("A218").Value = 11
$ws.Range("B218").Value = "Vega Monumental Concepción"
$ws.Range("C218").Value = "Bíobío"
$ws.Range("D218").Value = 44722
$ws.Range("E218").Value = 8
$ws.Range("F218").Value = 100112008
$ws.Range("G218").Value = "Coliflor"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 2200
$ws.Range("K218").Value = 650
$ws.Range("L218").Value = 700
$ws.Range("M218").Value = 677
$ws.Range("N218").Value = "$/unidad"
$ws.Range("O218").Value = "Región Metropolitana"
$ws.Range("P218").Value = 677
$ws.Range("Q218").Value = 1
$ws.Range("R218").Value = "Hortaliza"
